{"js": "// Add link to r-opensci blog post about R Research Compendiums\n// 1. Bump the Date paragraph from 2016-09-05 to 2016-09-06.\nconst body = context.document.body;\nconst allParas = body.paragraphs;\nallParas.load(\"items/style,items/text\");\nawait context.sync();\n\nfor (const p of allParas.items) {\n  if (p.style === \"Date\") {\n    p.insertText(\"2016-09-06\", Word.InsertLocation.replace);\n    break;\n  }\n}\nawait context.sync();\n\n// 2. Locate the \"Notes\" section hyperlink paragraph (styled \"First Paragraph\")\n//    that links to r-statistics.com; everything else is anchored relative to it.\nconst paras = body.paragraphs;\nparas.load(\"items/text,items/style\");\nawait context.sync();\n\nlet linkPara = null;\nfor (const p of paras.items) {\n  if (p.text.indexOf(\"r-statistics.com\") !== -1) {\n    linkPara = p;\n    break;\n  }\n}\n\n// 3. Insert a new explanatory paragraph right before it, in the same\n//    \"First Paragraph\" style used by the rest of the Notes section intro text.\nconst introPara = linkPara.insertParagraph(\n  \"These are links to resources that may be useful for writing this or as suggested resources in the final document that aren't easily printable for inclusion in the Dropbox folder.\",\n  Word.InsertLocation.before\n);\nintroPara.style = \"First Paragraph\";\n\n// 4. The original hyperlink paragraph is no longer the first paragraph of the\n//    section, so restyle it to \"Body Text\".\nlinkPara.style = \"Body Text\";\n\n// 5. Insert a new \"Body Text\" paragraph after it with a hyperlink to the\n//    rOpenSci rrrpkg repository.\nconst newLinkPara = linkPara.insertParagraph(\"https://github.com/ropensci/rrrpkg\", Word.InsertLocation.after);\nnewLinkPara.style = \"Body Text\";\nawait context.sync();\n\nconst newRange = newLinkPara.getRange(Word.RangeLocation.whole);\nnewRange.hyperlink = \"https://github.com/ropensci/rrrpkg\";\nawait context.sync();\n", "ps1": "# Apply the \"Add link to r-opensci blog post about R Research Compendiums\" edit.\n$d = $word.ActiveDocument\n\n# 1. Bump the Date paragraph from 2016-09-05 to 2016-09-06.\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Style.NameLocal -eq \"Date\") {\n        $d.Paragraphs($i).Range.Text = \"2016-09-06\"\n        break\n    }\n}\n\n# 2. Find the \"Notes\" hyperlink paragraph (style \"First Paragraph\") that links to\n#    r-statistics.com; everything else is anchored relative to it.\n$linkParaIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs($i).Range.Text -like \"*r-statistics.com*\") {\n        $linkParaIndex = $i\n        break\n    }\n}\n\n# 3. Insert a new explanatory paragraph right before it, in the same\n#    \"First Paragraph\" style used by the rest of the Notes section intro text.\n$linkPara = $d.Paragraphs($linkParaIndex)\n$linkPara.Range.InsertParagraphBefore()\n$introPara = $d.Paragraphs($linkParaIndex)\n$introPara.Style = \"First Paragraph\"\n$introPara.Range.Text = \"These are links to resources that may be useful for writing this or as suggested resources in the final document that aren't easily printable for inclusion in the Dropbox folder.\"\n\n# 4. The original link paragraph shifted down by one; restyle it to \"Body Text\"\n#    now that it is no longer the first paragraph of the section.\n$linkParaIndex = $linkParaIndex + 1\n$linkPara = $d.Paragraphs($linkParaIndex)\n$linkPara.Style = \"Body Text\"\n\n# 5. Add a new \"Body Text\" paragraph after it with a hyperlink to the rrrpkg repo.\n$linkPara.Range.InsertParagraphAfter()\n$newLinkParaIndex = $linkParaIndex + 1\n$newLinkPara = $d.Paragraphs($newLinkParaIndex)\n$newLinkPara.Style = \"Body Text\"\n$newLinkPara.Range.Text = \"https://github.com/ropensci/rrrpkg\"\n$rng = $d.Paragraphs($newLinkParaIndex).Range\n$rng.End = $rng.End - 1\n$d.Hyperlinks.Add($rng, \"https://github.com/ropensci/rrrpkg\", $null, $null, $rng.Text) | Out-Null\n"}
